$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.08312210549351147
$ws.Range("C2").Value = 0.6473643844303238
$ws.Range("D2").Value = 0.6561623056018049
$ws.Range("E2").Value = 0.8100384593349904
$ws.Range("F2").Value = 0.8340433064168259
$ws.Range("G2").Value = 15

$ws.Range("B3").Value = -0.01572770659491909
$ws.Range("C3").Value = 0.4445496780545906
$ws.Range("D3").Value = 0.3664663352916562
$ws.Range("E3").Value = 0.6053646300302457
$ws.Range("F3").Value = 0.6280045102432898
$ws.Range("G3").Value = 14

$ws.Range("B4").Value = -0.06939138655364743
$ws.Range("C4").Value = 0.5952380793966869
$ws.Range("D4").Value = 0.6460436807749044
$ws.Range("E4").Value = 0.803768424843191
$ws.Range("F4").Value = 0.8334651917952676
$ws.Range("G4").Value = 13

$ws.Range("B5").Value = 0.06184427096005268
$ws.Range("C5").Value = 0.33700112578256
$ws.Range("D5").Value = 0.209417853196624
$ws.Range("E5").Value = 0.4576219544521701
$ws.Range("F5").Value = 0.4735857100263252
$ws.Range("G5").Value = 12

$ws.Range("B6").Value = -0.06798794410060216
$ws.Range("C6").Value = 0.4968139668552824
$ws.Range("D6").Value = 0.377721656252307
$ws.Range("E6").Value = 0.6145906412013666
$ws.Range("F6").Value = 0.640631895303542
$ws.Range("G6").Value = 11

$ws.Range("B7").Value = 0.01694414359772521
$ws.Range("C7").Value = 0.2607354299704002
$ws.Range("D7").Value = 0.1060912636732037
$ws.Range("E7").Value = 0.3257165388389169
$ws.Range("F7").Value = 0.3428704965612517
$ws.Range("G7").Value = 10

$ws.Range("B8").Value = -0.06054918682133375
$ws.Range("C8").Value = 0.3330393258437492
$ws.Range("D8").Value = 0.1559239451819853
$ws.Range("E8").Value = 0.3948720617896198
$ws.Range("F8").Value = 0.4138719111052573
$ws.Range("G8").Value = 9

$ws.Range("B9").Value = -0.01715367655228583
$ws.Range("C9").Value = 0.2997021343775136
$ws.Range("D9").Value = 0.1700925356682539
$ws.Range("E9").Value = 0.4124227632760514
$ws.Range("F9").Value = 0.4405169522264152
$ws.Range("G9").Value = 8

$ws.Range("B10").Value = -0.05536040423422233
$ws.Range("C10").Value = 0.4661046776547557
$ws.Range("D10").Value = 0.3164365824987914
$ws.Range("E10").Value = 0.5625269615749909
$ws.Range("F10").Value = 0.6046490245027418
$ws.Range("G10").Value = 7

$ws.Range("B11").Value = -0.0196809274911882
$ws.Range("C11").Value = 0.3184982995745606
$ws.Range("D11").Value = 0.1151557242996088
$ws.Range("E11").Value = 0.3393460244346599
$ws.Range("F11").Value = 0.371109232533003
$ws.Range("G11").Value = 6

